# Updates crypto price/volume table cells to the latest scraped values
# (GitHub Actions scheduled refresh of cryptos.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.277.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.07%  "

$ws.Range("D3").Value = "'1.554.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.63%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'207.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.40%  "

$ws.Range("D7").Value = "'0.478"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.15%  "

$ws.Range("E8").Value = "  -1.70%  "

$ws.Range("E9").Value = "  -3.30%  "

$ws.Range("E10").Value = "  -4.61%  "

$ws.Range("E11").Value = "  -1.00%  "

$ws.Range("D12").Value = "'1.770.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.70%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.63%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.542.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.41%  "

$ws.Range("D16").Value = "'25.278.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.05%  "

$ws.Range("D17").Value = "'0.0₃0709"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.80%  "

$ws.Range("D18").Value = "'58.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.63%  "

$ws.Range("D19").Value = "'1.00"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'185.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.04%  "

$ws.Range("D22").Value = "'9.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.18%  "

$ws.Range("E23").Value = "  -3.82%  "

$ws.Range("E24").Value = "  -4.38%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "'140.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.88%  "

$ws.Range("E27").Value = "  -5.03%  "

$ws.Range("D28").Value = "'14.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.07%  "

$ws.Range("E29").Value = "  -5.16%  "

$ws.Range("E30").Value = "  -7.03%  "

$ws.Range("E31").Value = "  -3.40%  "

$ws.Range("E32").Value = "  -3.41%  "

$ws.Range("E33").Value = "  -5.14%  "

$ws.Range("E34").Value = "  -3.47%  "

$ws.Range("E35").Value = "  -3.54%  "

$ws.Range("D36").Value = "'1.085.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.00%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("E38").Value = "  -3.01%  "

$ws.Range("E39").Value = "  -4.99%  "

$ws.Range("D40").Value = "'0.766"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.32%  "

$ws.Range("E41").Value = "  -7.82%  "

$ws.Range("D42").Value = "'0.798"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.18%  "

$ws.Range("D43").Value = "'92.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.96%  "

$ws.Range("E44").Value = "  -1.67%  "

$ws.Range("D45").Value = "'1.684.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.67%  "

$ws.Range("E46").Value = "  -2.91%  "

$ws.Range("E47").Value = "  -2.13%  "

$ws.Range("E48").Value = "  -4.16%  "

$ws.Range("E49").Value = "  -4.09%  "

$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("E51").Value = "  -2.25%  "
